$wb = $excel.ActiveWorkbook

# Overview sheet: row 4 is for file 730c43ea-7b8b-436c-b1ff-838e0662ebf5.md
# Column G = "Latest HO Xliff Generate Date"
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G4").Value = "2016-08-25 20:44:58"

# zh-cn sheet: row 4 is for file 730c43ea-7b8b-436c-b1ff-838e0662ebf5.md
# Column H = "Correspond Handoff Datetime"
# Column K = "Correspond Handback DateTime"
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H4").Value = "2016-08-25 20:44:54"
$wsZhCn.Range("K4").Value = "2016-08-25 20:45:22"

# de-de sheet: row 4 is for file 730c43ea-7b8b-436c-b1ff-838e0662ebf5.md
# Column H = "Correspond Handoff Datetime"
# Column K = "Correspond Handback DateTime"
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H4").Value = "2016-08-25 20:44:58"
$wsDeDe.Range("K4").Value = "2016-08-25 20:45:29"
